# Trade #40 closed at 2026-02-17 13:27:15 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.69    # Current Capital
$summary.Range("B4").Value = -2.31      # Total P&L $
$summary.Range("B5").Value = -1.16      # Total P&L %
$summary.Range("B6").Value = 40         # Total Trades
$summary.Range("B7").Value = 16         # Winning Trades
$summary.Range("B9").Value = 40         # Win Rate %

# --- Strategy Status sheet (row 4 = MarketMaking) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.69       # Capital
$status.Range("D4").Value = 40          # Trades
$status.Range("E4").Value = -2.31       # P&L $
$status.Range("F4").Value = -2.31       # P&L %
$status.Range("G4").Value = 40          # Win Rate %

# --- New trade row (#40) appended to both "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(41, 1).Value = 40

    # Date/time columns must stay text (matching existing rows), not get
    # auto-coerced into date/time serial numbers.
    $ws.Cells.Item(41, 2).NumberFormat = "@"
    $ws.Cells.Item(41, 2).Value = "2026-02-17"
    $ws.Cells.Item(41, 2).ClearFormats()

    $ws.Cells.Item(41, 3).NumberFormat = "@"
    $ws.Cells.Item(41, 3).Value = "13:27:09"
    $ws.Cells.Item(41, 3).ClearFormats()

    $ws.Cells.Item(41, 4).Value = "MarketMaking"
    $ws.Cells.Item(41, 5).Value = "UP"
    $ws.Cells.Item(41, 6).Value = 0.72
    $ws.Cells.Item(41, 7).Value = 0.7481449999999999
    $ws.Cells.Item(41, 8).Value = "CLOSED"
    $ws.Cells.Item(41, 9).Value = 3.909
    $ws.Cells.Item(41, 10).Value = 0.03
    $ws.Cells.Item(41, 11).Value = 97.69
    $ws.Cells.Item(41, 12).Value = 0
    $ws.Cells.Item(41, 13).Value = 0
    $ws.Cells.Item(41, 14).Value = 0.6
    $ws.Cells.Item(41, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(41, 16).Value = "early_exit"
    $ws.Cells.Item(41, 17).Value = 0.15
}
